# Updated cryptos list on Sun Oct 22 11:14:05 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns on Sheet1 with the
# latest scraped values. Numeric-looking Price strings are written with a
# leading apostrophe so Excel stores them as literal text (matching the
# source data's text formatting, e.g. "28.56" / "0.0909") instead of
# auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.924.80'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').Value = '1.633.07'
$ws.Range('E3').Value = '  +1.90%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = '''214.44'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').Value = '''28.56'
$ws.Range('E8').Value = '  +0.76%  '
$ws.Range('E9').Value = '  +1.32%  '
$ws.Range('E10').Value = '  +0.76%  '
$ws.Range('D11').Value = '''0.0909'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = '1.867.64'
$ws.Range('E12').Value = '  +1.95%  '
$ws.Range('D13').Value = '1.629.27'
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').Value = '''0.563'
$ws.Range('E14').Value = '  +2.05%  '
$ws.Range('D15').Value = '''9.28'
$ws.Range('E15').Value = '  +16.35%  '
$ws.Range('D16').Value = '29.933.66'
$ws.Range('E16').Value = '  +0.74%  '
$ws.Range('E17').Value = '  +1.96%  '
$ws.Range('D18').Value = '''64.02'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = '''241.79'
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('E22').Value = '  +2.58%  '
$ws.Range('D23').Value = '''9.78'
$ws.Range('E23').Value = '  +3.67%  '
$ws.Range('E24').Value = '  +3.01%  '
$ws.Range('D25').Value = '''158.14'
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('D28').Value = '''6.61'
$ws.Range('E28').Value = '  +2.51%  '
$ws.Range('D30').Value = '''0.0487'
$ws.Range('E30').Value = '  +1.94%  '
$ws.Range('E31').Value = '  +4.17%  '
$ws.Range('D32').Value = '''3.37'
$ws.Range('E32').Value = '  +4.11%  '
$ws.Range('E33').Value = '  -0.40%  '
$ws.Range('D34').Value = '1.424.60'
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('E35').Value = '  +4.70%  '
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('E37').Value = '  -2.57%  '
$ws.Range('D38').Value = '''2.29'
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').Value = '''75.38'
$ws.Range('E40').Value = '  +12.52%  '
$ws.Range('E41').Value = '  +1.39%  '
$ws.Range('D42').Value = '''2.00'
$ws.Range('E42').Value = '  +3.03%  '
$ws.Range('D43').Value = '''0.828'
$ws.Range('E43').Value = '  +1.38%  '
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('E45').Value = '  +1.92%  '
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('D47').Value = '''52.70'
$ws.Range('E47').Value = '  -5.34%  '
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('D49').Value = '1.774.35'
$ws.Range('E49').Value = '  +1.96%  '
$ws.Range('D50').Value = '0.0₆0115'
$ws.Range('E50').Value = '  +9.59%  '
$ws.Range('D51').Value = '''90.32'
$ws.Range('E51').Value = '  +4.35%  '
